$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns
# to match the latest scrape. Price values that look like plain
# numbers are prefixed with a leading apostrophe so Excel keeps
# them as literal text (preserving formatting such as trailing
# zeros and the "thousands dot" notation) instead of silently
# converting them to floating point numbers.

$ws.Range("D2").Value = '28.900.46'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '1.886.32'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").Value = "'325.56"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = "'0.4591"
$ws.Range("E7").Value = '  +0.53%  '
$ws.Range("D8").Value = "'0.3891"
$ws.Range("E8").Value = '  +0.89%  '
$ws.Range("D9").Value = "'0.07863"
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = "'0.9869"
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = "'21.80"
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("D12").Value = '1.896.87'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").Value = "'7.023"
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").Value = "'5.678"
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = "'0.06932"
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").Value = "'88.28"
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = "'0.000009975"
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").Value = "'17.00"
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").Value = '28.908.78'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").Value = "'5.285"
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = "'10.97"
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = '2.124.41'
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").Value = "'2.084"
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").Value = "'155.17"
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").Value = "'19.28"
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").Value = "'5.976"
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("D29").Value = "'1.932"
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("D30").Value = "'117.47"
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("D31").Value = "'0.09337"
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").Value = "'0.9038"
$ws.Range("E32").Value = '  -1.17%  '
$ws.Range("D33").Value = "'5.281"
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("D34").Value = "'1.329"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = "'3.265"
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("D37").Value = "'0.05767"
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("D38").Value = "'0.02071"
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").Value = "'1.001"
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").Value = "'7.644"
$ws.Range("E40").Value = '  -0.64%  '
$ws.Range("D41").Value = "'0.5669"
$ws.Range("E41").Value = '  +1.15%  '
$ws.Range("D42").Value = "'0.1767"
$ws.Range("E42").Value = '  -0.59%  '
$ws.Range("D43").Value = "'9.693"
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").Value = "'2.262"
$ws.Range("E44").Value = '  +4.12%  '
$ws.Range("D45").Value = "'11.83"
$ws.Range("E45").Value = '  +1.82%  '
$ws.Range("D46").Value = "'0.5363"
$ws.Range("D47").Value = "'0.07033"
$ws.Range("E47").Value = '  -1.63%  '
$ws.Range("D48").Value = "'1.854"
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("D49").Value = "'2.560"
$ws.Range("E49").Value = '  +4.28%  '
$ws.Range("D50").Value = "'112.73"
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").Value = "'1.065"
$ws.Range("E51").Value = '  -4.82%  '
